# Updates the crypto price/volume table (and a few reordered rows) to match
# the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 27/28, 32/33 and 36/37 swapped position (B = Coin name, C = Link).
$rowData = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '27.694.19'; E = '  -0.06%  ' },
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '1.902.95'; E = '  +0.71%  ' },
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '0.9996'; E = '  -0.07%  ' },
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '312.29'; E = '  -0.31%  ' },
    @{ Row = 6; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '0.9993'; E = '  -0.14%  ' },
    @{ Row = 7; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '0.5228'; E = '  +7.69%  ' },
    @{ Row = 8; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.3777'; E = '  -0.31%  ' },
    @{ Row = 9; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.07242'; E = '  -1.25%  ' },
    @{ Row = 10; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '21.11'; E = '  +3.10%  ' },
    @{ Row = 11; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.8958'; E = '  -2.55%  ' },
    @{ Row = 12; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.07626'; E = '  -0.74%  ' },
    @{ Row = 13; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.896.71'; E = '  +0.29%  ' },
    @{ Row = 14; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.444'; E = '  -0.33%  ' },
    @{ Row = 15; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '92.00'; E = '  +1.22%  ' },
    @{ Row = 16; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.000'; E = '  -0.09%  ' },
    @{ Row = 17; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.000008718'; E = '  -0.90%  ' },
    @{ Row = 18; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '0.9997'; E = '  -0.11%  ' },
    @{ Row = 19; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '27.711.01'; E = '  -0.09%  ' },
    @{ Row = 20; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '14.45'; E = '  -0.51%  ' },
    @{ Row = 21; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '5.127'; E = '  +0.17%  ' },
    @{ Row = 22; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.146.68'; E = '  -0.08%  ' },
    @{ Row = 23; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '10.80'; E = '  -0.03%  ' },
    @{ Row = 24; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '6.574'; E = '  -0.36%  ' },
    @{ Row = 25; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '153.28'; E = '  +0.02%  ' },
    @{ Row = 26; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.869'; E = '  -2.22%  ' },
    @{ Row = 27; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '2.164'; E = '  +2.37%  ' },
    @{ Row = 28; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '18.27'; E = '  -0.56%  ' },
    @{ Row = 29; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '114.50'; E = '  -1.12%  ' },
    @{ Row = 30; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.833'; E = '  -1.31%  ' },
    @{ Row = 31; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.08986'; E = '  +0.62%  ' },
    @{ Row = 32; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.853'; E = '  +4.83%  ' },
    @{ Row = 33; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '3.171'; E = '  +0.63%  ' },
    @{ Row = 34; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.236'; E = '  +1.18%  ' },
    @{ Row = 35; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.7688'; E = '  +0.99%  ' },
    @{ Row = 36; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '2.616'; E = '  +2.77%  ' },
    @{ Row = 37; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.02075'; E = '  +2.01%  ' },
    @{ Row = 38; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '3.054'; E = '  +2.80%  ' },
    @{ Row = 39; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '1.089'; E = '  -0.15%  ' },
    @{ Row = 40; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.5484'; E = '  +0.58%  ' },
    @{ Row = 41; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.05271'; E = '  +0.41%  ' },
    @{ Row = 42; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '6.644'; E = '  -4.22%  ' },
    @{ Row = 43; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '113.14'; E = '  +3.03%  ' },
    @{ Row = 44; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '8.456'; E = '  +1.63%  ' },
    @{ Row = 45; B = 'Algorand'; C = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D = '0.1507'; E = '  -0.76%  ' },
    @{ Row = 46; B = 'Decentraland'; C = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D = '0.4778'; E = '  +0.01%  ' },
    @{ Row = 47; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '10.39'; E = '  -1.64%  ' },
    @{ Row = 48; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '0.9992'; E = '  -0.16%  ' },
    @{ Row = 49; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '1.612'; E = '  -1.31%  ' },
    @{ Row = 50; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '66.36'; E = '  -1.65%  ' },
    @{ Row = 51; B = 'Cronos'; C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = '0.05998'; E = '  -0.93%  ' }
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C

    # The Price column holds text that looks numeric (e.g. "92.00", "1.000").
    # Assigning it directly would make Excel coerce it to a Double and drop the
    # formatting (trailing zeros, thousand-separator dots). Prefix with a literal
    # apostrophe to force text, then restore the cell's original style so no
    # stray quote-prefix formatting is left behind.
    $dCell = $ws.Range("D$r")
    $dStyle = $dCell.Style
    $dCell.Value = "`'" + $item.D
    $dCell.Style = $dStyle

    $ws.Range("E$r").Value = $item.E
}